# Fruta / hortaliza, semanal
# Insert a new weekly price record for "Feria Lagunitas de Puerto Montt - Poroto verde"
# as row 65, pushing the existing rows 65-83 down to 66-84.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 65, shifting rows 65-83 down to 66-84.
$ws.Rows.Item(65).Insert()

# Populate the newly inserted row 65 with the new weekly data point.
$ws.Cells.Item(65, 1).Value = 4
$ws.Cells.Item(65, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(65, 3).Value = "Los Lagos"
$ws.Cells.Item(65, 4).Value = 44726
$ws.Cells.Item(65, 5).Value = 10
$ws.Cells.Item(65, 6).Value = 100112031
$ws.Cells.Item(65, 7).Value = "Poroto verde"
$ws.Cells.Item(65, 8).Value = "Magnum"
$ws.Cells.Item(65, 9).Value = "Primera"
$ws.Cells.Item(65, 10).Value = 40
$ws.Cells.Item(65, 11).Value = 26000
$ws.Cells.Item(65, 12).Value = 26000
$ws.Cells.Item(65, 13).Value = 26000
$ws.Cells.Item(65, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(65, 15).Value = "Perú"
$ws.Cells.Item(65, 16).Value = 1040
$ws.Cells.Item(65, 17).Value = 25
$ws.Cells.Item(65, 18).Value = "Hortaliza"
